$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F11").Value = 7
$ws.Range("G11").Value = 1019.13
$ws.Range("B14").Value = 8798.6
$ws.Range("F64").Value = 37
$ws.Range("G64").Value = 365.19
$ws.Range("B71").Value = 54593.08
$ws.Range("F93").Value = 60
$ws.Range("G93").Value = 2388
$ws.Range("B94").Value = 24384.08
$ws.Range("F105").Value = 220
$ws.Range("G105").Value = 17857.4
$ws.Range("F107").Value = 131
$ws.Range("G107").Value = 5881.9
$ws.Range("F116").Value = 12
$ws.Range("G116").Value = 1619.4
$ws.Range("F123").Value = 6
$ws.Range("G123").Value = 280.44
$ws.Range("F127").Value = 107
$ws.Range("G127").Value = 5435.6
$ws.Range("F130").Value = 59
$ws.Range("G130").Value = 11900.89
$ws.Range("B143").Value = 277423
$ws.Range("F192").Value = 44
$ws.Range("G192").Value = 2177.12
$ws.Range("B205").Value = 26628.05
$ws.Range("F235").Value = 24
$ws.Range("G235").Value = 2026.08
$ws.Range("F238").Value = 16
$ws.Range("G238").Value = 811.84
$ws.Range("F249").Value = 88
$ws.Range("G249").Value = 3516.48
$ws.Range("B250").Value = 27075.21
$ws.Range("F296").Value = 69
$ws.Range("G296").Value = 5853.27
$ws.Range("B301").Value = 13071.94
$ws.Range("F311").Value = 43
$ws.Range("G311").Value = 6170.5
$ws.Range("B319").Value = 30348.02
$ws.Range("F335").Value = 4
$ws.Range("G335").Value = 80.84
$ws.Range("F346").Value = 74
$ws.Range("G346").Value = 4437.78
$ws.Range("F353").Value = 1
$ws.Range("G353").Value = 32.07
$ws.Range("F360").Value = 51
$ws.Range("G360").Value = 2483.19
$ws.Range("B382").Value = 125766.89
$ws.Range("F387").Value = 238
$ws.Range("G387").Value = 24444.98
$ws.Range("F425").Value = 0
$ws.Range("G425").Value = 0
$ws.Range("F434").Value = 170
$ws.Range("G434").Value = 17187
$ws.Range("F442").Value = 155
$ws.Range("G442").Value = 9165.15
$ws.Range("F456").Value = 16
$ws.Range("G456").Value = 2418.72
$ws.Range("B467").Value = 409695.03
$ws.Range("F475").Value = 0
$ws.Range("G475").Value = 0
$ws.Range("F481").Value = 12
$ws.Range("G481").Value = 2618.76
$ws.Range("B482").Value = 35344.39
$ws.Range("F485").Value = 0
$ws.Range("G485").Value = 0
$ws.Range("B487").Value = 13018.95
$ws.Range("F492").Value = 2
$ws.Range("G492").Value = 72.2
$ws.Range("F496").Value = 14
$ws.Range("G496").Value = 285.46
$ws.Range("B500").Value = 1276.97
$ws.Range("F524").Value = 52
$ws.Range("G524").Value = 50430.12
$ws.Range("B525").Value = 50430.12
$ws.Range("F536").Value = 94
$ws.Range("G536").Value = 9080.4
$ws.Range("F538").Value = 78
$ws.Range("G538").Value = 1918.02
$ws.Range("B542").Value = 55176.48
$ws.Range("F573").Value = 32
$ws.Range("G573").Value = 412.8
$ws.Range("B582").Value = 1268.95
$ws.Range("F598").Value = 1035
$ws.Range("G598").Value = 6810.3
$ws.Range("B605").Value = 110661.74
$ws.Range("F616").Value = 2
$ws.Range("G616").Value = 69.62
$ws.Range("F617").Value = 62
$ws.Range("G617").Value = 2158.22
$ws.Range("B620").Value = 10592.72
$ws.Range("F669").Value = 683
$ws.Range("G669").Value = 4644.4
$ws.Range("F670").Value = 438
$ws.Range("G670").Value = 3000.3
$ws.Range("F674").Value = 324
$ws.Range("G674").Value = 5355.72
$ws.Range("F675").Value = 148
$ws.Range("G675").Value = 3904.24
$ws.Range("B677").Value = 42977.84
$ws.Range("F712").Value = 46
$ws.Range("G712").Value = 4531
$ws.Range("B718").Value = 69898.89
$ws.Range("F747").Value = 34
$ws.Range("G747").Value = 1707.82
$ws.Range("F754").Value = 143
$ws.Range("G754").Value = 3809.52
$ws.Range("B761").Value = 54039.17
$ws.Range("F768").Value = 101
$ws.Range("G768").Value = 2747.2
$ws.Range("B770").Value = 68759.12
$ws.Range("F795").Value = 66
$ws.Range("G795").Value = 2185.26
$ws.Range("F803").Value = 28
$ws.Range("G803").Value = 1225.28
$ws.Range("B804").Value = 36341.34
$ws.Range("F833").Value = 8
$ws.Range("G833").Value = 734
$ws.Range("B839").Value = 3996.82
$ws.Range("F895").Value = 0
$ws.Range("G895").Value = 0
$ws.Range("B901").Value = 50354.95
$ws.Range("F921").Value = 78
$ws.Range("G921").Value = 6130.02
$ws.Range("B933").Value = 36847.78
$ws.Range("F936").Value = 65
$ws.Range("G936").Value = 2431
$ws.Range("F941").Value = 123
$ws.Range("G941").Value = 4600.2
$ws.Range("B942").Value = 13505.1
$ws.Range("F999").Value = 1397
$ws.Range("G999").Value = 227864.67
$ws.Range("F1001").Value = 231
$ws.Range("G1001").Value = 17819.34
$ws.Range("B1005").Value = 268669.59
$ws.Range("B1012").Value = 2586347.36
$ws.Range("B1013").Value = 2586347.36
